$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.750.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.034.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.99%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'228.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.54%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.05%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'60.24"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.19%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.46%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0820"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.05%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.59%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'14.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.81%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.335.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.94%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'20.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.46%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.770"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.70%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.60%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.041.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.55%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.702.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.91%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'69.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.51%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'5.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -7.08%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0822"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'222.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.25%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.27%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.89%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'167.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.128"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.96%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'18.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.12%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.30%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.45%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +8.85%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.19%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.72%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.27%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.64%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.82%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.12%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'17.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.29%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.533.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.96%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0216"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.36%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'96.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.47%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.03%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0911"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.25%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.68%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.94%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.67%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.06%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'7.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.57%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.224.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.97%  "
$ws.Range("E51").Style = "Normal"
